# "contingencies with rene fine"
# Adds a tiny 2x2 "disconnected_elements" flag table to Sheet1:
#   B1 = 0   (bold, thin box border, centered/top aligned)
#   A2 = 0   (bold, thin box border, centered/top aligned)
#   B2 = "disconnected_elements" (plain, stored as a shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values first.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the formatted style once on B1: bold font, thin box border,
# centered horizontally / top vertically.
$cell1 = $ws.Range("B1")
$cell1.Font.Bold = $true
$cell1.HorizontalAlignment = -4108   # xlCenter
$cell1.VerticalAlignment = -4160     # xlTop
$cell1.Borders.LineStyle = 1         # xlContinuous
$cell1.Borders.Weight = 2            # xlThin

# Clone that exact formatting onto A2 via copy/paste-special so both
# cells end up sharing a single cell style (matches how Excel itself
# would dedupe the style table instead of minting a second near-
# identical xf entry).
$cell1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
